# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 09:29"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6001017
$ws.Range("C4").Value = 652
$ws.Range("D4").Value = 3314305
$ws.Range("E4").Value = 2503049
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 183663

# Row 6 - India
$ws.Range("B6").Value = 3314953
$ws.Range("C6").Value = 7204
$ws.Range("D6").Value = 2524518
$ws.Range("E6").Value = 729783
$ws.Range("G6").Value = 23
$ws.Range("H6").Value = 60652

# Row 58 - Armenia
$ws.Range("B58").Value = 43270
$ws.Range("C58").Value = 203
$ws.Range("D58").Value = 36988
$ws.Range("E58").Value = 5418
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 864

# Row 63 - Afganistan
$ws.Range("B63").Value = 38129
$ws.Range("C63").Value = 16
$ws.Range("D63").Value = 29046
$ws.Range("E63").Value = 7682

# Row 109 - Hungria
$ws.Range("B109").Value = 5379
$ws.Range("C109").Value = 91
$ws.Range("D109").Value = 3757
$ws.Range("E109").Value = 1008

# Row 151 - Georgia
$ws.Range("B151").Value = 1447
$ws.Range("C151").Value = 11
$ws.Range("D151").Value = 1190
$ws.Range("E151").Value = 238
